$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'332.45"
$ws.Range("E2").Value = "'1.25%"
$ws.Range("D3").Value = "'45.87"
$ws.Range("E3").Value = "'4.68%"
$ws.Range("D4").Value = "'5.643"
$ws.Range("E4").Value = "'2.49%"
$ws.Range("E5").Value = "'4.32%"
$ws.Range("D6").Value = "'2.067"
$ws.Range("E6").Value = "'4.14%"
$ws.Range("E7").Value = "'3.64%"
$ws.Range("D8").Value = "'0.9889"
$ws.Range("E8").Value = "'4.18%"
$ws.Range("D9").Value = "'2.583"
$ws.Range("E9").Value = "'-1.16%"
$ws.Range("E10").Value = "'3.18%"
$ws.Range("D11").Value = "'0.1931"
$ws.Range("E11").Value = "'4.24%"
$ws.Range("D12").Value = "'10.40"
$ws.Range("E12").Value = "'-2.38%"
$ws.Range("D13").Value = "'0.09968"
$ws.Range("E13").Value = "'1.39%"
$ws.Range("D14").Value = "'0.04675"
$ws.Range("E14").Value = "'0.22%"
$ws.Range("E15").Value = "'-0.43%"
$ws.Range("D16").Value = "'0.001276"
$ws.Range("E16").Value = "'1.25%"
$ws.Range("D17").Value = "'0.006128"
$ws.Range("E17").Value = "'3.54%"
$ws.Range("E18").Value = "'0.62%"
$ws.Range("D19").Value = "'0.3366"
$ws.Range("E19").Value = "'-3.17%"
$ws.Range("D20").Value = "'0.1403"
$ws.Range("E20").Value = "'-0.19%"
$ws.Range("D21").Value = "'0.2656"
$ws.Range("E21").Value = "'4.46%"
$ws.Range("D22").Value = "'0.04211"
$ws.Range("E22").Value = "'3.73%"
$ws.Range("D23").Value = "'0.001316"
$ws.Range("E23").Value = "'4.90%"
$ws.Range("D24").Value = "'0.004649"
$ws.Range("E24").Value = "'7.20%"
$ws.Range("D25").Value = "'0.0001284"
$ws.Range("E25").Value = "'7.12%"
$ws.Range("D26").Value = "'0.0003751"
$ws.Range("D38").Value = "'0.02783"
$ws.Range("E38").Value = "'7.68%"
$ws.Range("D39").Value = "'0.05740"
$ws.Range("E39").Value = "'1.30%"
$ws.Range("D40").Value = "'0.007774"
$ws.Range("E40").Value = "'3.21%"
$ws.Range("D41").Value = "'0.1436"
$ws.Range("E41").Value = "'2.92%"
$ws.Range("D42").Value = "'0.007288"
$ws.Range("E42").Value = "'-3.31%"
$ws.Range("D43").Value = "'0.002120"
$ws.Range("E43").Value = "'5.36%"
$ws.Range("D44").Value = "'0.009068"
$ws.Range("E44").Value = "'8.20%"
$ws.Range("D45").Value = "'0.3409"
$ws.Range("D46").Value = "'0.00007112"
$ws.Range("E46").Value = "'-0.14%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.44%"
$ws.Range("D48").Value = "'0.0005817"
$ws.Range("E48").Value = "'0.10%"
$ws.Range("B49").Value = "'CoinbaseStockToken"
$ws.Range("C49").Value = "'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "'0.003508"
$ws.Range("E49").Value = "'-0.60%"
$ws.Range("B50").Value = "'BOLO"
$ws.Range("C50").Value = "'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D50").Value = "'0.003499"
$ws.Range("E50").Value = "'-2.04%"
$ws.Range("D51").Value = "'0.00002106"
$ws.Range("E51").Value = "'0.44%"
